$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed crypto price/volume data.
# Numeric-looking "Price" text values are protected from Excel's automatic
# number conversion by briefly switching to a text number format, then
# restoring the original cell style so no formatting changes leak into the file.

$ws.Range("D2").Value = "30.171.59"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "1.856.77"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  +0.43%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.92"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -1.57%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +0.33%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4767"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -3.07%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2813"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -4.52%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06507"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -3.28%  "
$ws.Range("D10").Value = "1.870.04"
$ws.Range("E10").Value = "  -0.99%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07347"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -0.14%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.37"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -4.63%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.143"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -0.33%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.11"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -1.32%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6447"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").Value = "30.147.07"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.20"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -2.09%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007604"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("D20").Value = "2.127.34"
$ws.Range("E20").Value = "  -0.54%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.21%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.274"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -0.18%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "216.04"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +12.77%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.092"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -2.13%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.301"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -2.42%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.72"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +2.10%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.42"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.11%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.903"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -2.38%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.428"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -3.27%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.234"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -4.79%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09137"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.83%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.967"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -4.36%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05023"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -4.28%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7422"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -0.21%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +1.77%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.698"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -0.59%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01830"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -0.49%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.609"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("B39").Value = "PaxosStandard"
$ws.Range("C39").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.110"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +10.91%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9086"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.047"
$ws.Range("D41").Style = $style
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.99"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.898"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -1.08%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4247"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.393"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -2.66%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1307"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -5.86%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.561"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +9.13%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.11"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -10.06%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.882"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.32"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -2.41%  "
